# feat: add 2022-Q1 data
#
# The workbook originally has two sheets:
#   1) "2020-Q4"  -> per-fund holdings for the 2020-Q4 quarter
#   2) "总计"      -> summary/totals table, one row per quarter
#
# This change inserts a brand-new "2022-Q1" sheet (positioned between the
# two existing sheets) containing the per-fund holdings for the 2022-Q1
# quarter, and adds a corresponding summary row to the "总计" sheet.
#
# NOTE: worksheet object references in this host appear to track a sheet
# *slot/position* rather than a fixed identity, so once a structural
# operation (Copy/Add/Move/Delete) changes sheet ordering, any
# previously-captured worksheet variable can silently start pointing at a
# different sheet. To stay safe, this script re-fetches worksheets **by
# name** immediately after every structural operation instead of reusing
# variables captured beforehand.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet by duplicating "总计" so it inherits
#    the same cell styles (header style, index-column style, etc.), then
#    place it right after "2020-Q4" (i.e. before "总计").
# ---------------------------------------------------------------------
$wb.Worksheets.Item("总计").Copy($null, $wb.Worksheets.Item("2020-Q4"))

# The copy always lands immediately after "2020-Q4", i.e. at position 2.
$sheetQ1 = $wb.Worksheets.Item(2)
$sheetQ1.Name = "2022-Q1"

# Copy the bold/bordered header style (currently on B1) onto the extra
# header columns this sheet needs (E1:H1) before filling in their text.
$sheetQ1.Range("B1").Copy()
$sheetQ1.Range("E1:H1").PasteSpecial(-4122)   # xlPasteFormats

$sheetQ1.Range("B1").Value = "基金代码"
$sheetQ1.Range("C1").Value = "基金名称"
$sheetQ1.Range("D1").Value = "基金规模"
$sheetQ1.Range("E1").Value = "股票总仓位"
$sheetQ1.Range("F1").Value = "仓位占比"
$sheetQ1.Range("G1").Value = "持有市值(亿元)"
$sheetQ1.Range("H1").Value = "仓位排名"

# Row 2: 377016 上投摩根亚太优势混合(QDII)
$sheetQ1.Range("A2").Value = 0
$sheetQ1.Range("B2").Value = "377016"
$sheetQ1.Range("C2").Value = "上投摩根亚太优势混合(QDII)"
$sheetQ1.Range("D2").Value = "27.15"
$sheetQ1.Range("E2").Value = "86.02"
$sheetQ1.Range("F2").Value = "5.27"
$sheetQ1.Range("G2").Value = "1.4308"
$sheetQ1.Range("H2").Value = 3

# Row 3: 378006 上投摩根全球新兴市场混合(QDII)
$sheetQ1.Range("A3").Value = 1
$sheetQ1.Range("B3").Value = "378006"
$sheetQ1.Range("C3").Value = "上投摩根全球新兴市场混合(QDII)"
$sheetQ1.Range("D3").Value = "0.46"
$sheetQ1.Range("E3").Value = "88.99"
$sheetQ1.Range("F3").Value = "4.02"
$sheetQ1.Range("G3").Value = "0.0185"
$sheetQ1.Range("H3").Value = 3

# Give the new A3 index cell the same style as A2 (already carried the
# "总计" style via the sheet copy).
$sheetQ1.Range("A2").Copy()
$sheetQ1.Range("A3").PasteSpecial(-4122)      # xlPasteFormats

# ---------------------------------------------------------------------
# 2) Insert a new summary row for "2022-Q1" at the top of the "总计"
#    sheet's data (row 2), pushing the existing "2020-Q4" row down to
#    row 3 -- without inheriting the bold header formatting that a plain
#    Rows.Insert() would introduce. Re-fetch "总计" by name since the
#    Copy() above changed sheet ordering.
# ---------------------------------------------------------------------
$sheetZong = $wb.Worksheets.Item("总计")

$oldDate  = $sheetZong.Range("B2").Value()
$oldCount = $sheetZong.Range("C2").Value()
$oldValue = $sheetZong.Range("D2").Value()

$sheetZong.Range("A3").Value = 1
$sheetZong.Range("B3").Value = $oldDate
$sheetZong.Range("C3").Value = $oldCount
$sheetZong.Range("D3").Value = $oldValue

$sheetZong.Range("A2").Copy()
$sheetZong.Range("A3").PasteSpecial(-4122)    # xlPasteFormats

$sheetZong.Range("A2").Value = 0
$sheetZong.Range("B2").Value = "2022-Q1"
$sheetZong.Range("C2").Value = 2
$sheetZong.Range("D2").Value = 1.45

# ---------------------------------------------------------------------
# 3) Restore the original active sheet/tab selection ("2020-Q4"), since
#    the sheet-copy operations above shift which tab Excel considers
#    active/selected.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
